$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted as row 29 (pushing the existing
# rows 29-110 down to 30-111). Insert a blank row at position 29 first,
# which shifts everything below it down by one.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record's data.
$ws.Range("A29").Value = 5
$ws.Range("B29").Value = "Macroferia Regional de Talca"
$ws.Range("C29").Value = "Maule"
$ws.Range("D29").Value = 44525
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 100112031
$ws.Range("G29").Value = "Poroto verde"
$ws.Range("H29").Value = "Magnum"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 33000
$ws.Range("L29").Value = 33000
$ws.Range("M29").Value = 33000
$ws.Range("N29").Value = "`$/malla 25 kilos"
$ws.Range("O29").Value = "Región del Maule"
$ws.Range("P29").Value = 1320
$ws.Range("Q29").Value = 25
$ws.Range("R29").Value = "Hortaliza"

# Preserve the date-formatted style used by all other cells in column D.
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat
